# "Finished revising following committee comments sent back to Curt"
#
# This edit turns the single "Sheet1" table into two tabs:
#   - the original data, renamed "pre Dec2015"
#   - a new, revised copy named "Rev Dec2015" with updated numbers,
#     carrying the committee's requested changes.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the original sheet -----------------------------------------
$orig = $wb.Worksheets.Item(1)
$orig.Name = "pre Dec2015"

# --- 2. Duplicate it (keeps styles/column widths/shared formatting intact)
#        and place the copy right after the original ------------------------
$orig.Copy($null, $orig) | Out-Null
$rev = $wb.ActiveSheet
$rev.Name = "Rev Dec2015"

# --- 3. Update the revised sheet's content ---------------------------------
# (values are written in the same order the original author entered them,
# so new shared-string entries land in the same slots)

# Row 9: new grid-cell-mean drifter speed ranges
$rev.Range("C9").Value = "1-19"
$rev.Range("B9").Value = "1-20"
$rev.Range("D9").Value = "1-36"
# ... relabelled as the grid-cell-mean drifter speed range
$rev.Range("A9").Value = "Drifters Speed range (mean in grid cells) (cm s-1)"

# Row 10: individual-drifter mean speed label + revised values
$rev.Range("A10").Value = "Drifters (individual) Mean Speed (cm s-1)"

# Row 11: individual-drifter standard-deviation label + revised values
$rev.Range("A11").Value = "Drifters (individual) Standard Deviation (cm s-1)"

# Row 6: ADCP Speed values (cm s-1) for the three ADCP stations
$rev.Range("B6").Value = "14.9, 6.8, 0.4"
$rev.Range("C6").Value = "12.5, 3.7, 0.7"
$rev.Range("D6").Value = "21.5, 11, 1.2"

# Rows 7-8: replace the AVERAGE/STDEV formulas with the committee-revised
# static mean / standard-deviation figures
$rev.Range("B7").Value = 7.4
$rev.Range("C7").Value = 5.6
$rev.Range("D7").Value = 11.2

$rev.Range("B8").Value = 7.3
$rev.Range("C8").Value = 6.1
$rev.Range("D8").Value = 10.1

# Row 10/11 revised numeric values
$rev.Range("B10").Value = 8
$rev.Range("C10").Value = 7.1
$rev.Range("D10").Value = 12.3

$rev.Range("B11").Value = 6.5
$rev.Range("C11").Value = 5.8
$rev.Range("D11").Value = 8.1

# Row 12: revised drifter counts, now formatted with a thousands separator
$rev.Range("B12").Value = 1314
$rev.Range("C12").Value = 1580
$rev.Range("D12").Value = 2461
$rev.Range("B12:D12").NumberFormat = "#,##0"

# --- 4. Selection / active-tab bookkeeping ----------------------------------
$rev.Range("D6").Select() | Out-Null
